$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename "Student Name" -> "Student_Name"
$ws.Range("A1").Value = "Student_Name"
$ws.Range("B1").Value = "Marks"
$ws.Range("C1").Value = "Progress"

# Existing student row (Rajesh) gains a rank value; marks column loses its
# date-style numeric formatting (now plain/General)
$ws.Range("A2").Value = "Rajesh"
$ws.Range("B2").Value = 80
$ws.Range("B2").NumberFormat = "General"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 45

# New student row
$ws.Range("A3").Value = "veera"
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 34
$ws.Range("D3").Value = 67

# New "rank" header column, added last in the shared-string table
$ws.Range("D1").Value = "rank"

$ws.Range("D3").Select()
